# Append a new attendance record (row 3) below the existing header (row 1)
# and first data row (row 2) of the trainer-attendance sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the new row's cells to Text first so values such as the date-like
# "2025-03-23" and time-like "07:29:30" strings are stored verbatim instead
# of being auto-converted to Excel date/time serial numbers.
$ws.Range("A3:F3").NumberFormat = "@"

$ws.Range("A3").Value = "2025-03-23"
$ws.Range("B3").Value = "BILAL"
$ws.Range("C3").Value = "Bilal"
$ws.Range("D3").Value = "07:29:30"
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = "manager"

# Drop the temporary number-format override so the new row ends up
# unstyled, matching the look of the existing data row (row 2).
$ws.Range("A3:F3").ClearFormats()
